$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A44").Value = "GRT-USD"
